$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that currently sits in the middle of
#    the "... ukoliko postoji popust, nakon cega se brise." sentence.
# ---------------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# ---------------------------------------------------------------------------
# 2) Re-create "_GoBack" as a collapsed bookmark right after the closing
#    parenthesis of "... 10% za obicne korisnike)" (this is where Word last
#    left the cursor / the new edit point).
#
#    Bookmarks.Add on a zero-length Range that sits exactly at a paragraph
#    boundary gets mis-positioned by this host, so the bookmark is built by
#    temporarily inserting a one-character placeholder, bookmarking that
#    single character, and then deleting the placeholder text again - the
#    bookmark collapses to an empty span in place, exactly like Word does
#    when the bookmarked text is removed.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("za obične korisnike)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endPos = $rng.End

$insPos = $d.Range($endPos, $endPos)
$insPos.InsertAfter("X")
$bmRange = $d.Range($endPos, $endPos + 1)
$d.Bookmarks.Add("_GoBack", $bmRange)
$placeholder = $d.Range($endPos, $endPos + 1)
$placeholder.Text = ""

# ---------------------------------------------------------------------------
# 3) The "_Hlk164719481" bookmark (around "(20% za privilegovane i 10)")
#    keeps its name/position - only its numeric id shifts from 0 to 1, which
#    happens automatically now that "_GoBack" occupies id 0.
# ---------------------------------------------------------------------------
